$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("?" bucket) count drops by one: 488 -> 487
$ws.Range("B2").Value = 487

# Rows 10-12 get new "value" (place-name) / date-range data while keeping
# their original "count" (column B), except B11 which becomes 14.
# Row 13 (Viena) is untouched.
$ws.Range("A10").Value = "Landsberg"
$ws.Range("A11").Value = "Macau"
$ws.Range("B11").Value = 14
$ws.Range("A12").Value = "Avignon"

# The C:H columns hold numeric-looking / date-looking values that must stay
# plain text (matching the workbook's original inlineStr cell type). Typing
# them straight into .Value lets Excel auto-convert to numbers/dates (and
# picking up a new number format if forced to text). Routing them through a
# "=""literal""" formula and then Copy/PasteSpecial-as-values keeps them as
# plain strings with no formula residue and no number-format change.
$ws.Range("C10").Formula = '="16230729"'
$ws.Range("D10").Formula = '="17571009"'
$ws.Range("E10").Formula = '="1623-07-29"'
$ws.Range("F10").Formula = '="1757-10-09"'
$ws.Range("G10").Formula = '="1623"'
$ws.Range("H10").Formula = '="1757"'

$ws.Range("C11").Formula = '="16280000"'
$ws.Range("D11").Formula = '="17490201"'
$ws.Range("E11").Formula = '="1628"'
$ws.Range("F11").Formula = '="1749-02-01"'
$ws.Range("G11").Formula = '="1628"'
$ws.Range("H11").Formula = '="1749"'

$ws.Range("C12").Formula = '="16150926"'
$ws.Range("D12").Formula = '="17370927"'
$ws.Range("E12").Formula = '="1615-09-26"'
$ws.Range("F12").Formula = '="1737-09-27"'
$ws.Range("G12").Formula = '="1615"'
$ws.Range("H12").Formula = '="1737"'

$ws.Range("C10:H12").Copy()
$ws.Range("C10:H12").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = 0
